$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 28 (item 24) with the new "Oil Engine + Filter" service entry
$ws.Range("B28").Value = "Oil Engine + Filter"
$ws.Range("C28").Value = "15W-40"
$ws.Range("D28").Value = "9 ltr"
$ws.Range("E28").Value = 301308
$ws.Range("F28").Formula = "=7000+E28"
$ws.Range("G28").Value = "25/03/2024"

# Move the active selection to F27, matching the author's final cursor position
[void]$ws.Range("F27").Select()
